# Apply "4 times the wind power" re-run results to res_bus/p_mw sheet.
# Updates the bus active power (p_mw) results: column B (bus p_mw) and
# column I (another bus's recomputed p_mw, which is exactly 4x the
# original value) for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 98.33082189154038
$ws.Range("I2").Value = -157.0811126327369

$ws.Range("B3").Value = 117.0706584358068
$ws.Range("I3").Value = -159.0775654985071

$ws.Range("B4").Value = 117.325339161175
$ws.Range("I4").Value = -156.8510428125631

$ws.Range("B5").Value = 115.0606543929562
$ws.Range("I5").Value = -151.0364367016142

$ws.Range("B6").Value = 113.9344408466623
$ws.Range("I6").Value = -146.7455716947669

$ws.Range("B7").Value = 124.9166644827703
$ws.Range("I7").Value = -146.7455716947669

$ws.Range("B8").Value = 127.0055097470711
$ws.Range("I8").Value = -145.9459847787531

$ws.Range("B9").Value = 118.2961956206209
$ws.Range("I9").Value = -144.7164313135622

$ws.Range("B10").Value = 91.28743648603427
$ws.Range("I10").Value = -148.7483106211975

$ws.Range("B11").Value = 60.87609920852447
$ws.Range("I11").Value = -151.5632085849628

$ws.Range("B12").Value = 48.41386117890852
$ws.Range("I12").Value = -150.3273690591114

$ws.Range("B13").Value = 44.14817850620693
$ws.Range("I13").Value = -148.949464562333

$ws.Range("B14").Value = 46.35620485809864
$ws.Range("I14").Value = -146.9668410300159

$ws.Range("B15").Value = 48.12812539728429
$ws.Range("I15").Value = -144.0224502166446

$ws.Range("B16").Value = 51.82599569214176
$ws.Range("I16").Value = -146.9002087870148

$ws.Range("B17").Value = 56.03530794402536
$ws.Range("I17").Value = -150.6366432436073

$ws.Range("B18").Value = 62.82017645987594
$ws.Range("I18").Value = -150.2758233616954

$ws.Range("B19").Value = 67.23338340370105
$ws.Range("I19").Value = -151.9479154973845

$ws.Range("B20").Value = 59.47376180224251
$ws.Range("I20").Value = -151.6197831309072

$ws.Range("B21").Value = 42.2969152050116
$ws.Range("I21").Value = -153.5106301775812

$ws.Range("B22").Value = 33.89216090204673
$ws.Range("I22").Value = -157.7763509417866

$ws.Range("B23").Value = 29.21873965769328
$ws.Range("I23").Value = -162.3777473452619

$ws.Range("B24").Value = 28.99733288624952
$ws.Range("I24").Value = -165.6150685854118

$ws.Range("B25").Value = 49.33643119517637
$ws.Range("I25").Value = -168
